$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 20: a new scoring result row ---
# B20 picks up the same numeric style as the rest of column B (style index 4,
# the "0.00000_ " custom number format already used throughout the sheet).
$ws.Range("B20").Value = 0.94538679999999997

# C20/D20 in this sheet were entered without the column's default numeric
# style (same as C19/D19 immediately above them) - copy that "no special
# format" look from row 19 before writing the values so no new style gets
# minted.
$ws.Range("C19:D19").Copy()
$ws.Range("C20:D20").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C20").Value = 0.92373000000000005
$ws.Range("D20").Value = 0.89822999999999997

# Label the new row: headline+snippet+abstract whole-text emotion feature.
$ws.Range("F20").Value = "headline+snippet+abstract整段文字的emotion"

# --- Update F19's label to be more specific (headline+snippet only) ---
$ws.Range("F19").Value = "headline+snippet整段文字的emotion"

# Leave the selection on the newly added cell, like the author did.
$ws.Range("C20").Select()
